$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 10708.272
$ws.Cells.Item(2, 9).Value = 20048.2
$ws.Cells.Item(2, 10).Value = 2925
$ws.Cells.Item(2, 11).Value = 20048.2
$ws.Cells.Item(2, 12).Value = 2925
$ws.Cells.Item(2, 13).Value = -19935.2
$ws.Cells.Item(2, 14).Value = -3151

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 289.03226
$ws.Cells.Item(18, 9).Value = 289.03226
$ws.Cells.Item(18, 11).Value = 289.03226
$ws.Cells.Item(18, 13).Value = -5.032260000000008

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 6901
$ws.Cells.Item(70, 9).Value = 3626.25
$ws.Cells.Item(70, 10).Value = 20000
$ws.Cells.Item(70, 11).Value = 10878.75
$ws.Cells.Item(70, 12).Value = 60000
$ws.Cells.Item(70, 13).Value = -10608.75
$ws.Cells.Item(70, 14).Value = -60540

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 6901
$ws.Cells.Item(73, 9).Value = 3626.25
$ws.Cells.Item(73, 10).Value = 20000
$ws.Cells.Item(73, 11).Value = 10878.75
$ws.Cells.Item(73, 12).Value = 60000
$ws.Cells.Item(73, 13).Value = -9942.75
$ws.Cells.Item(73, 14).Value = -61872

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1068.3903
$ws.Cells.Item(61, 9).Value = 1046.4517
$ws.Cells.Item(61, 10).Value = 1136.4
$ws.Cells.Item(61, 11).Value = 1046.4517
$ws.Cells.Item(61, 12).Value = 1136.4
$ws.Cells.Item(61, 13).Value = -834.4517000000001
$ws.Cells.Item(61, 14).Value = -1560.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 858.6389
$ws.Cells.Item(97, 9).Value = 676.931
$ws.Cells.Item(97, 10).Value = 1611.4286
$ws.Cells.Item(97, 11).Value = 676.931
$ws.Cells.Item(97, 12).Value = 1611.4286
$ws.Cells.Item(97, 13).Value = -180.931
$ws.Cells.Item(97, 14).Value = -2603.4286

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1063.1578
$ws.Cells.Item(102, 9).Value = 698.2308
$ws.Cells.Item(102, 10).Value = 1853.8334
$ws.Cells.Item(102, 11).Value = 698.2308
$ws.Cells.Item(102, 12).Value = 1853.8334
$ws.Cells.Item(102, 13).Value = 923.7692
$ws.Cells.Item(102, 14).Value = -5097.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 22530070
$ws.Cells.Item(132, 9).Value = 30304696
$ws.Cells.Item(132, 10).Value = 4204165
$ws.Cells.Item(132, 11).Value = 90914088
$ws.Cells.Item(132, 12).Value = 12612495
$ws.Cells.Item(132, 13).Value = -90911558
$ws.Cells.Item(132, 14).Value = -12617555

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1068.3903
$ws.Cells.Item(136, 9).Value = 1046.4517
$ws.Cells.Item(136, 10).Value = 1136.4
$ws.Cells.Item(136, 11).Value = 3139.3551
$ws.Cells.Item(136, 12).Value = 3409.2
$ws.Cells.Item(136, 13).Value = -589.3551000000002
$ws.Cells.Item(136, 14).Value = -8509.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 43479936
$ws.Cells.Item(105, 9).Value = 1645.2106
$ws.Cells.Item(105, 10).Value = 250001820
$ws.Cells.Item(105, 11).Value = 1645.2106
$ws.Cells.Item(105, 12).Value = 250001820
$ws.Cells.Item(105, 13).Value = 101.7893999999999
$ws.Cells.Item(105, 14).Value = -250005314

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2927750
$ws.Cells.Item(134, 9).Value = 874.37036
$ws.Cells.Item(134, 11).Value = 2623.11108
$ws.Cells.Item(134, 13).Value = -88.11108000000013

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 14252.533
$ws.Cells.Item(105, 9).Value = 24437.125
$ws.Cells.Item(105, 10).Value = 2613
$ws.Cells.Item(105, 11).Value = 24437.125
$ws.Cells.Item(105, 12).Value = 2613
$ws.Cells.Item(105, 13).Value = -22690.125
$ws.Cells.Item(105, 14).Value = -6107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 66668930
$ws.Cells.Item(132, 9).Value = 2400
$ws.Cells.Item(132, 11).Value = 7200
$ws.Cells.Item(132, 13).Value = -4670

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 13158835
$ws.Cells.Item(134, 9).Value = 1051.931
$ws.Cells.Item(134, 11).Value = 3155.793
$ws.Cells.Item(134, 13).Value = -620.7930000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 24451.355
$ws.Cells.Item(12, 9).Value = 2.8333333
$ws.Cells.Item(12, 11).Value = 8.499999900000001
$ws.Cells.Item(12, 13).Value = 164.5000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(99, 8).Value = 1250.25
$ws.Cells.Item(99, 9).Value = 600.4
$ws.Cells.Item(99, 10).Value = 2333.3333
$ws.Cells.Item(99, 11).Value = 1801.2
$ws.Cells.Item(99, 12).Value = 6999.999899999999
$ws.Cells.Item(99, 13).Value = 444.8000000000002
$ws.Cells.Item(99, 14).Value = -11491.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 887.61
$ws.Cells.Item(131, 10).Value = 930.57776
$ws.Cells.Item(131, 12).Value = 2791.73328
$ws.Cells.Item(131, 14).Value = -12871.73328

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 6666
$ws.Cells.Item(62, 10).Value = 6666
$ws.Cells.Item(62, 12).Value = 6666
$ws.Cells.Item(62, 14).Value = -8038

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(65, 8).Value = 6666
$ws.Cells.Item(65, 10).Value = 6666
$ws.Cells.Item(65, 12).Value = 19998
$ws.Cells.Item(65, 14).Value = -26862

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 11872.579
$ws.Cells.Item(132, 9).Value = 7624
$ws.Cells.Item(132, 10).Value = 27804.75
$ws.Cells.Item(132, 11).Value = 22872
$ws.Cells.Item(132, 12).Value = 83414.25
$ws.Cells.Item(132, 13).Value = -20342
$ws.Cells.Item(132, 14).Value = -88474.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(63, 8).Value = 11139
$ws.Cells.Item(63, 10).Value = 11139
$ws.Cells.Item(63, 12).Value = 11139
$ws.Cells.Item(63, 14).Value = -12637

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(66, 8).Value = 11139
$ws.Cells.Item(66, 10).Value = 11139
$ws.Cells.Item(66, 12).Value = 33417
$ws.Cells.Item(66, 14).Value = -40905

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1580.6666
$ws.Cells.Item(82, 9).Value = 1571
$ws.Cells.Item(82, 10).Value = 1600
$ws.Cells.Item(82, 11).Value = 1571
$ws.Cells.Item(82, 12).Value = 1600
$ws.Cells.Item(82, 13).Value = -1210
$ws.Cells.Item(82, 14).Value = -2322

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1580.6666
$ws.Cells.Item(85, 9).Value = 1571
$ws.Cells.Item(85, 10).Value = 1600
$ws.Cells.Item(85, 11).Value = 1571
$ws.Cells.Item(85, 12).Value = 1600
$ws.Cells.Item(85, 13).Value = -323
$ws.Cells.Item(85, 14).Value = -4096

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1129.8572
$ws.Cells.Item(93, 9).Value = 976.375
$ws.Cells.Item(93, 10).Value = 1334.5
$ws.Cells.Item(93, 11).Value = 976.375
$ws.Cells.Item(93, 12).Value = 1334.5
$ws.Cells.Item(93, 13).Value = 271.625
$ws.Cells.Item(93, 14).Value = -3830.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 23260404
$ws.Cells.Item(132, 9).Value = 35716230
$ws.Cells.Item(132, 10).Value = 9530.666999999999
$ws.Cells.Item(132, 11).Value = 107148690
$ws.Cells.Item(132, 12).Value = 28592.001
$ws.Cells.Item(132, 13).Value = -107146160
$ws.Cells.Item(132, 14).Value = -33652.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 36006436
$ws.Cells.Item(136, 9).Value = 13229240
$ws.Cells.Item(136, 11).Value = 39687720
$ws.Cells.Item(136, 13).Value = -39685170

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1197.909
$ws.Cells.Item(81, 9).Value = 1117.5
$ws.Cells.Item(81, 11).Value = 2235
$ws.Cells.Item(81, 13).Value = -1174

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 1197.909
$ws.Cells.Item(84, 9).Value = 1117.5
$ws.Cells.Item(84, 11).Value = 11175
$ws.Cells.Item(84, 13).Value = -5871

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 47789.6
$ws.Cells.Item(132, 9).Value = 116279.555
$ws.Cells.Item(132, 10).Value = 9264
$ws.Cells.Item(132, 11).Value = 348838.665
$ws.Cells.Item(132, 12).Value = 27792
$ws.Cells.Item(132, 13).Value = -346308.665
$ws.Cells.Item(132, 14).Value = -32852
